{"js": "// Replace each \"a\u00f7b=\" division prompt in the document with its new value.\n// Several prompts repeat (e.g. \"77\u00f75=\" and \"56\u00f74=\" each occur twice) so we\n// search for every distinct old string once and then walk the matches in\n// document order, applying the replacements for that string in sequence.\nconst replacements = {\n  \"91\u00f77=\": [\"92\u00f73=\"],\n  \"32\u00f78=\": [\"92\u00f72=\"],\n  \"72\u00f79=\": [\"64\u00f75=\"],\n  \"41\u00f73=\": [\"76\u00f74=\"],\n  \"77\u00f75=\": [\"45\u00f72=\", \"25\u00f75=\"],\n  \"20\u00f75=\": [\"31\u00f77=\"],\n  \"86\u00f72=\": [\"35\u00f78=\"],\n  \"68\u00f77=\": [\"93\u00f75=\"],\n  \"95\u00f78=\": [\"83\u00f77=\"],\n  \"56\u00f74=\": [\"58\u00f74=\", \"87\u00f79=\"],\n  \"98\u00f79=\": [\"65\u00f72=\"],\n  \"88\u00f78=\": [\"92\u00f78=\"],\n  \"55\u00f73=\": [\"47\u00f75=\"],\n  \"19\u00f77=\": [\"83\u00f79=\"],\n  \"51\u00f73=\": [\"65\u00f73=\"],\n  \"60\u00f77=\": [\"73\u00f76=\"],\n  \"17\u00f76=\": [\"19\u00f79=\"],\n  \"41\u00f75=\": [\"25\u00f73=\"],\n  \"29\u00f72=\": [\"21\u00f79=\"],\n  \"74\u00f74=\": [\"47\u00f79=\"],\n  \"75\u00f76=\": [\"72\u00f75=\"],\n  \"28\u00f77=\": [\"60\u00f73=\"],\n  \"16\u00f73=\": [\"15\u00f74=\"],\n};\n\nconst body = context.document.body;\n\nfor (const oldText of Object.keys(replacements)) {\n  const newTexts = replacements[oldText];\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  const count = Math.min(newTexts.length, results.items.length);\n  for (let i = 0; i < count; i++) {\n    results.items[i].insertText(newTexts[i], \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"a\u00f7b=\" division prompt with its new value, walking every\n# table cell in document order. Several prompts repeat (e.g. \"77\u00f75=\" and\n# \"56\u00f74=\" each occur twice), so matching is done positionally against the\n# ordered list below rather than by a plain global find/replace.\n$oldValues = @(\n  \"91\u00f77=\", \"32\u00f78=\", \"72\u00f79=\", \"41\u00f73=\", \"77\u00f75=\",\n  \"20\u00f75=\", \"86\u00f72=\", \"68\u00f77=\", \"95\u00f78=\", \"77\u00f75=\",\n  \"56\u00f74=\", \"98\u00f79=\", \"88\u00f78=\", \"55\u00f73=\", \"19\u00f77=\",\n  \"51\u00f73=\", \"60\u00f77=\", \"17\u00f76=\", \"41\u00f75=\", \"29\u00f72=\",\n  \"74\u00f74=\", \"56\u00f74=\", \"75\u00f76=\", \"28\u00f77=\", \"16\u00f73=\"\n)\n$newValues = @(\n  \"92\u00f73=\", \"92\u00f72=\", \"64\u00f75=\", \"76\u00f74=\", \"45\u00f72=\",\n  \"31\u00f77=\", \"35\u00f78=\", \"93\u00f75=\", \"83\u00f77=\", \"25\u00f75=\",\n  \"58\u00f74=\", \"65\u00f72=\", \"92\u00f78=\", \"47\u00f75=\", \"83\u00f79=\",\n  \"65\u00f73=\", \"73\u00f76=\", \"19\u00f79=\", \"25\u00f73=\", \"21\u00f79=\",\n  \"47\u00f79=\", \"87\u00f79=\", \"72\u00f75=\", \"60\u00f73=\", \"15\u00f74=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n$next = 0\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($next -ge $oldValues.Count) { continue }\n    $cellRange = $t.Cell($r, $c).Range\n    $cellText = $cellRange.Text\n    $expected = $oldValues[$next]\n    if ($cellText -eq ($expected + \"`r`a\")) {\n      $cellRange.Text = $newValues[$next]\n      $next = $next + 1\n    }\n  }\n}\n"}
